# Generate Report for Handback
# Update the timestamp strings recorded in the handback status report.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 17:07:24"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-22 17:07:19"
$wsZhCn.Range("K2").Value = "2016-08-22 17:07:46"

# "de-de" sheet: Correspond Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-08-22 17:07:53"
